$d = $word.ActiveDocument

# Each of these text spans was previously split across two (or three)
# adjacent runs purely because of a trailing/leading space run with
# identical formatting. Re-"typing" the same text via Find/Replace
# coalesces the runs Word normally keeps merged.
#
# Replace = 1 (wdReplaceOne) is used deliberately instead of wdReplaceAll:
# several of these needles (e.g. "bold ") also occur as a substring inside
# later, differently-formatted text ("bold italics"), and a replace-all
# pass would needlessly re-touch/rewrite that unrelated run.

$d.Content.Find.Execute("Regular text ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Regular text ", 1)

$d.Content.Find.Execute("bold ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "bold ", 1)

$d.Content.Find.Execute("This is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This is ", 1)

$d.Content.Find.Execute(", and this is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", and this is ", 1)

$d.Content.Find.Execute("Some people use ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Some people use ", 1)

$d.Content.Find.Execute("single underlines for ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "single underlines for ", 1)

$d.Content.Find.Execute("Above the line is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Above the line is ", 1)

$d.Content.Find.Execute(" and below the line is ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " and below the line is ", 1)
